# Apply the "investment" (事業投資) sheet update described in commit
# "#5: insurance, claim, debt, investment done".
#
# Sheet6 ("事業投資") previously only had a raw data dump (row 1 was an
# accidental duplicate of row 2, and only columns A-G were populated).
# This change turns row 1 into proper column headers and appends the
# standard legislator/property metadata columns (H:N) that the other
# sheets in this workbook already have, and fixes two amount cells that
# were stored as text instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# --- 1. Build header row (row 1) --------------------------------------
# Give the new header cells (H1:N1) the same bold/border style as the
# existing header cells by copying the formatting across first.
$ws.Range("B1:G1").Copy($ws.Range("H1"))
$ws.Range("B1").Copy($ws.Range("N1"))

$ws.Range("B1").Value = "owner"
$ws.Range("C1").Value = "company"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- 2. Fix numeric amount cells that were stored as text --------------
$ws.Range("E3").Value = 2000000
$ws.Range("E4").Value = 2500000

# --- 3. Fill in the new metadata columns (H:N) for each data row -------
# Column J holds an ISO-formatted date ("2012-04-30"); entering that as a
# plain value would make Excel auto-convert it into a date serial number.
# Instead, enter it as a text formula (="2012-04-30") and immediately
# collapse the formula down to its static text result with a values-only
# paste, which keeps it a genuine text cell without disturbing styles.
$rows = 2, 3, 4
foreach ($r in $rows) {
    $ws.Range("H$r").Value = "investment"
    $ws.Range("I$r").Value = "normal"

    $ws.Range("J$r").Formula = "=""2012-04-30"""
    $ws.Range("J$r").Copy()
    $ws.Range("J$r").PasteSpecial(-4163)

    $ws.Range("K$r").Value = "楊麗環"
    $ws.Range("L$r").Value = 960
    $ws.Range("M$r").Value = "tmp700a1"
}

$ws.Range("N2").Value = 144
$ws.Range("N3").Value = 145
$ws.Range("N4").Value = 146

Write-Host "sheet6 (investment) headers and metadata columns added"
